$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.531.80"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.793.16"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'329.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").Value = "'0.4399"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("D8").Value = "'0.3736"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.67%  "
$ws.Range("D9").Value = "'45.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").Value = "'0.07590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'1.133"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "'22.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'1.004"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "'6.218"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'7.490"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").Value = "1.794.91"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'0.00001088"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'0.06704"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("D19").Value = "'80.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'17.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").Value = "'6.218"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.62%  "
$ws.Range("D23").Value = "28.543.76"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'11.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "'2.444"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").Value = "'20.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "'152.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("D28").Value = "'2.332"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").Value = "2.002.78"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").Value = "'1.306"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "'130.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "'3.981"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("D33").Value = "'5.785"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "'0.09250"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").Value = "'0.2237"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("D36").Value = "'12.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'0.06260"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'5.184"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'0.02315"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").Value = "'0.6577"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("D41").Value = "'1.194"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.56%  "
$ws.Range("D42").Value = "'1.424"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.99%  "
$ws.Range("D43").Value = "'7.981"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'13.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").Value = "'0.6077"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'3.812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "'127.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").Value = "'2.009"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'0.07001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "'1.135"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.30%  "
